# New reading pulled in from Adafruit IO: append as the next row after the
# existing data, mirroring the layout/values of the most recent readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 69

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"

# "25" looks numeric, so a plain Value assignment would be auto-typed as a
# number by Excel. Every other cell in this sheet stores its value as text
# (inline/shared strings), so force this one to stay text too: format as
# Text, assign, then drop the formatting again so no visible style sticks
# to the cell (matches the rest of the sheet, which uses the default style).
$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "25"
$ws.Range("C$newRow").ClearFormats()

$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
